$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("main refs")

# Correct a data entry error before the column shift: row 13's "Group" value
# (column D, "pts") should actually be "Western".
$ws.Range("D13").Value = "Western"

# Remove the redundant/unused "Species-group" column (column B); every
# value in it was the constant "PLV" and is no longer needed.
$ws.Columns.Item(2).Delete()

$ws.Select()
$ws.Range("A1:J13").Select()
